$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing values that changed (International = B, National = C)
$ws.Range("C3").Value = 2578876
$ws.Range("C4").Value = 1380710
$ws.Range("B9").Value = 532317
$ws.Range("C10").Value = 6874502
$ws.Range("B16").Value = 363420
$ws.Range("B17").Value = 632122
$ws.Range("B18").Value = 972858
$ws.Range("C18").Value = 6056063
$ws.Range("B19").Value = 1948049
$ws.Range("B20").Value = 1117360
$ws.Range("C20").Value = 4725820

# Add new row 21 with Provisional Occupancy Forecast data
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = 43544
$ws.Range("B21").Value = 667329
$ws.Range("C21").Value = 1729850
$ws.Range("D21").Formula = '=B21+C21/Hoja2!$A$2'

# Extend table range to include new row
$ws.ListObjects.Item("Tabla1").Resize($ws.Range("A1:D21"))

# Update selection to match target state
$ws.Range("C6").Select()
